# Conversions work, formatting tidied up a bit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out the leftover "Variations" column (C): header + the few
# scratch notes that were only half-filled in.
$ws.Range("C1").ClearContents()
$ws.Range("C30").ClearContents()
$ws.Range("C37").ClearContents()
$ws.Range("C43").ClearContents()

# Tidy up the selection/view left over from editing.
$ws.Range("D15").Select()
